$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============================================================
# Step 1: copy cell formats to their new positions. Order matters:
# always read a format from a range before that range becomes a
# paste destination itself.
# ============================================================

# Good Drivers data-row format: old row 15 -> all new data rows 17:27
$ws.Range("A15:E15").Copy()
$ws.Range("A17:E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Good Drivers" title format: old row 13 -> new row 15
$ws.Range("A13").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Good Drivers column-header format: old row 14 -> new row 16
$ws.Range("A14:E14").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Totals row format: old row 7 (A:C) -> new row 9
$ws.Range("A7:C7").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Bad Drivers data-row format: row 6 -> new rows 7:8
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ============================================================
# Step 2: clear cells that become obsolete/blank in the new layout
# ============================================================
$ws.Range("A13:E14").Clear()
$ws.Range("B15:E15").Clear()

# ============================================================
# Step 3: write the Bad Drivers table content (rows 3-8)
# ============================================================
$ws.Range("A3").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.800"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 34
$ws.Range("D3").Value = 89.9
$ws.Range("A4").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.897"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 352
$ws.Range("D4").Value = 94
$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.60.1.2"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = 97.4
$ws.Range("A6").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.824"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 29
$ws.Range("D6").Value = 97.8
$ws.Range("A7").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.908"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 35
$ws.Range("D7").Value = 98.6
$ws.Range("A8").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.40.0.4"
$ws.Range("B8").Value = 36
$ws.Range("C8").Value = 506
$ws.Range("D8").Value = 98.6

# ============================================================
# Step 4: write the Totals row (row 9)
# ============================================================
$ws.Range("A9").Value = "Totals:"
$ws.Range("B9").Value = 51
$ws.Range("C9").Value = 985

# ============================================================
# Step 5: write section titles/headers
# ============================================================
$ws.Range("A15").Value = "Good Drivers (Roaming > 99.8%)"
$ws.Range("A16").Value = "Adapter-Driver"
$ws.Range("B16").Value = "Total Samples"
$ws.Range("D16").Value = "Good Roaming Calculation (%)"
$ws.Range("E16").Value = "Driver Vintage"

# ============================================================
# Step 6: write the Good Drivers table content (rows 17-27).
# The "Driver Vintage" column (E) holds YYYY-MM-DD text that Excel
# would otherwise auto-convert to a date serial number, so for each
# of those cells we temporarily force text format, assign the
# string, then restore the original (General, s=4-style) format by
# re-pasting formats from column D of the same row.
# ============================================================
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B17").Value = 449371
$ws.Range("D17").Value = 99.9
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2024-11-10"
$ws.Range("D17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A18").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.4.0.1088"
$ws.Range("B18").Value = 86281
$ws.Range("D18").Value = 99.9
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2024-08-07"
$ws.Range("D18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A19").Value = "MediaTek Wi-Fi 6E MT7922 (RZ616) 160MHz PCIe Adapter - 3.3.0.1030"
$ws.Range("B19").Value = 17891
$ws.Range("D19").Value = 100
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2024-05-09"
$ws.Range("D19").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.3.1"
$ws.Range("B20").Value = 10661
$ws.Range("D20").Value = 100
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2022-08-29"
$ws.Range("D20").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.150.0.3"
$ws.Range("B21").Value = 14239
$ws.Range("D21").Value = 100
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2022-05-23"
$ws.Range("D21").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.100.1.1"
$ws.Range("B22").Value = 265400
$ws.Range("D22").Value = 99.9
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2022-05-01"
$ws.Range("D22").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B23").Value = 77999
$ws.Range("D23").Value = 99.9
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2021-08-18"
$ws.Range("D23").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B24").Value = 34244
$ws.Range("D24").Value = 100
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2021-04-27"
$ws.Range("D24").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B25").Value = 59673
$ws.Range("D25").Value = 100
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2020-08-05"
$ws.Range("D25").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B26").Value = 113652
$ws.Range("D26").Value = 100
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "2020-01-06"
$ws.Range("D26").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A27").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B27").Value = 56018
$ws.Range("D27").Value = 100
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2019-12-14"
$ws.Range("D27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "edit complete"
